# Update cell values to reflect the new step data / selector documentation,
# mirroring the commit "agregando variable p para guardar datos entre paginas".

$wb = $excel.ActiveWorkbook

# --- Sheet "ProductSearchAdd" ---
$ws1 = $wb.Worksheets.Item("ProductSearchAdd")

# D8: "Datos" column - tidy up blank lines and swap the example product.
$ws1.Range("D8").Value = "producto:case para iphone`nproducto: tenis adidas`nnumeroProducto: 1`nelegir segundo producto de la lista`nProducto: nintendo`nNombre articulo`nMarca`nprecio"

# E8: "Observaciones" column - append the new botonCarrito selector used to
# pass the cart button reference ("p") between pages. This cell keeps its
# original "quote prefix" cell format (style index 6), so re-enter it with a
# leading apostrophe - same as Excel's own text-entry convention - to keep
# that format after the rewrite instead of falling back to the plain style.
$ws1.Range("E8").Value = "'SearchPage: `nsearchBar; //css=`"[class*=searchBarContainer] input`"`nsearchIcon; //css=`"button[data-automation-id='search-icon']`"`n`nResultsPage:`nlistaArticulos; //[class*='col']`nnombreArticulo; //css=`"[class*=product_name] p div`"`nmarcaArticulo;  //css=`"[class*=product_brand] a`"`nprecioArticulo; //css=`"[class*=product_price_] p`"`n`n`nProductPage:`nnombreArticulo;//css=`"[data-automation-id='product-name'] div`"`nmarcaArticulo; //css=`"[data-automation-id='brand'] a`"`nprecioArticulo; //css=`"[data-automation-id='list-price'] span`"`n`naddToCart; //css=`"[data-automation-id='add-button']`"`nproductAddedValidation; //css=`"[data-automation-id='add-button']`"`nbotonCarrito: //css=`"[data-automation-id='add-button']`""

# --- Sheet "CartReviewDelete" ---
$ws2 = $wb.Worksheets.Item("CartReviewDelete")

# E7: "Observaciones" column - add a verification selector for added products.
$ws2.Range("E7").Value = "CartPage:`nverificarProdAgregados; //css=`"[data-automation-id='added-to-cart']`"`nbotonCarrito: //css=`"[data-automation-id='go-to-cart']`"`nlistaArticulos;//`nmosaicoArticulo;//`nnombreArticulo;//`nprecioArticulo;//"

# The added line makes the cell wrap to one more row, so the row grows taller.
$ws2.Rows.Item(7).RowHeight = 100.8

# Move the active selection to D10 on CartReviewDelete, matching where the
# author ended up editing, but keep ProductSearchAdd as the active tab (it
# was the active sheet before/after this edit too).
$ws2.Activate()
$ws2.Range("D10").Select()
$ws1.Activate()
